$d = $word.ActiveDocument

# Paragraph 1 is the title ("On Pilgrimage - September 1965", style Heading1)
# Paragraph 2 is the byline ("By Dorothy Day", bold run, no style / Normal)
# Replace both paragraphs together with:
#   - a Title-styled paragraph whose text is split word-by-word/space into
#     separate runs
#   - an Authors-styled paragraph with "Dorothy Day" similarly split into
#     separate runs (no bold)
$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)
$r = $d.Range($p1.Range.Start, $p2.Range.End)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="Title"/></w:pPr>
<w:r><w:t xml:space="preserve">On</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">Pilgrimage</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">-</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">September</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">1965</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="Authors"/></w:pPr>
<w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">Day</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r.InsertXML($xml)
